$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "cube" placeholder type key/value rows appended below the existing
# goal_error_* rows (row 12 was the last used row).
$ws.Range("A13").Value = "cube_field"
$ws.Range("B13").Value = "Field"

$ws.Range("A14").Value = "cube_pond"
$ws.Range("B14").Value = "Pond"

$ws.Range("A15").Value = "cube_house_1"
$ws.Range("A16").Value = "cube_house_2"
$ws.Range("A17").Value = "cube_house_3"

$ws.Range("B15").Value = "House 01"
$ws.Range("B16").Value = "House 02"
$ws.Range("B17").Value = "House 03"

# Move the active selection to where the user's cursor ended up after
# entering the new rows.
[void]$ws.Range("B18").Select()
